$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dict_detail import-template header row was:
#   A1=dict_id  B1=lbl  C1=val  D1=is_locked  E1=is_enabled  F1=order_by  G1=rem
# Drop the "is_locked" / "is_enabled" header columns (D1:E1) entirely so the
# remaining columns (order_by, rem, ...) shift left and the now-unused
# shared-string entries for those two columns are dropped as well.
$ws.Range("D1:E1").EntireColumn.Delete()
